$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add partial hexadecimal color codes next to the matching RGB rows.
$ws.Range("G10").Value = "#404040"
$ws.Range("G11").Value = "#254061"
$ws.Range("G4").Value = "#FFCC0"
$ws.Range("G9").Value = "#98487"

# Reflect the final active selection left by the edit.
$ws.Range("G9").Select()
